# NB train and test fixed (they were being trained with the labels)
# Adds 7 new worksheets (ClinicalNB, ClinicalNB1, ClinicalNB2, LabNB,
# ClinicalNB3, LabNB1, LabNB2) at the end of the workbook, each reporting
# Naive-Bayes results for Clinical / Lab datasets.

$wb = $excel.ActiveWorkbook
$origActiveSheet = $wb.ActiveSheet

function Add-SheetAtEnd($name) {
    $count = $wb.Worksheets.Count
    $lastSheet = $wb.Worksheets.Item($count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $name
    return $ws
}

function Style-HeaderLikeCell($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

# ---------------------------------------------------------------------------
# Sheet: ClinicalNB  (Conf_id / Dataset / Base table, 3 data rows)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "ClinicalNB"

$ws.Range("B1").Value = "Conf_id"
$ws.Range("C1").Value = "Dataset"
$ws.Range("D1").Value = "Base"
$ws.Range("E1").Value = "P1"
$ws.Range("F1").Value = "P2"
$ws.Range("G1").Value = "P3"
$ws.Range("H1").Value = "P4"
$ws.Range("I1").Value = "P5"
$ws.Range("J1").Value = "Promedio"
Style-HeaderLikeCell $ws.Range("B1:J1")
Style-HeaderLikeCell $ws.Range("A2:A4")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Clinical"
$ws.Range("D2").Value = "Gaussian"
$ws.Range("E2").Value = 0.7931283430792027
$ws.Range("F2").Value = 0.7920333509430529
$ws.Range("G2").Value = 0.7955894503171277
$ws.Range("H2").Value = 0.8154063829277552
$ws.Range("I2").Value = 0.7845617124982224
$ws.Range("J2").Value = 0.7961438479530721

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Clinical"
$ws.Range("D3").Value = "Categorical"
$ws.Range("E3").Value = 0.9216641870017656
$ws.Range("F3").Value = 0.9200612080359398
$ws.Range("G3").Value = 0.9114584800680022
$ws.Range("H3").Value = 0.9086402948228408
$ws.Range("I3").Value = 0.9205627034707913
$ws.Range("J3").Value = 0.916477374679868

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Clinical"
$ws.Range("D4").Value = "Bernoulli"
$ws.Range("E4").Value = 0.9216641870017656
$ws.Range("F4").Value = 0.9200612080359398
$ws.Range("G4").Value = 0.9114584800680022
$ws.Range("H4").Value = 0.9086402948228408
$ws.Range("I4").Value = 0.9205627034707913
$ws.Range("J4").Value = 0.916477374679868

# ---------------------------------------------------------------------------
# Sheet: ClinicalNB1  (Base / P1..P5 / Promedio / Accuracy, 1 data row)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "ClinicalNB1"

$ws.Range("B1").Value = "Base"
$ws.Range("C1").Value = "P1"
$ws.Range("D1").Value = "P2"
$ws.Range("E1").Value = "P3"
$ws.Range("F1").Value = "P4"
$ws.Range("G1").Value = "P5"
$ws.Range("H1").Value = "Promedio"
$ws.Range("I1").Value = "Accuracy"
Style-HeaderLikeCell $ws.Range("B1:I1")
Style-HeaderLikeCell $ws.Range("A2")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Categorical"
$ws.Range("C2").Value = 0.9216641870017656
$ws.Range("D2").Value = 0.9200612080359398
$ws.Range("E2").Value = 0.9114584800680022
$ws.Range("F2").Value = 0.9086402948228408
$ws.Range("G2").Value = 0.9205627034707913
$ws.Range("H2").Value = 0.916477374679868
$ws.Range("I2").Value = 0.9334

# ---------------------------------------------------------------------------
# Sheet: ClinicalNB2  (identical data to ClinicalNB1)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "ClinicalNB2"

$ws.Range("B1").Value = "Base"
$ws.Range("C1").Value = "P1"
$ws.Range("D1").Value = "P2"
$ws.Range("E1").Value = "P3"
$ws.Range("F1").Value = "P4"
$ws.Range("G1").Value = "P5"
$ws.Range("H1").Value = "Promedio"
$ws.Range("I1").Value = "Accuracy"
Style-HeaderLikeCell $ws.Range("B1:I1")
Style-HeaderLikeCell $ws.Range("A2")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Categorical"
$ws.Range("C2").Value = 0.9216641870017656
$ws.Range("D2").Value = 0.9200612080359398
$ws.Range("E2").Value = 0.9114584800680022
$ws.Range("F2").Value = 0.9086402948228408
$ws.Range("G2").Value = 0.9205627034707913
$ws.Range("H2").Value = 0.916477374679868
$ws.Range("I2").Value = 0.9334

# ---------------------------------------------------------------------------
# Sheet: LabNB  (Conf_id / Dataset / Base table, 3 data rows)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "LabNB"

$ws.Range("B1").Value = "Conf_id"
$ws.Range("C1").Value = "Dataset"
$ws.Range("D1").Value = "Base"
$ws.Range("E1").Value = "P1"
$ws.Range("F1").Value = "P2"
$ws.Range("G1").Value = "P3"
$ws.Range("H1").Value = "P4"
$ws.Range("I1").Value = "P5"
$ws.Range("J1").Value = "Promedio"
Style-HeaderLikeCell $ws.Range("B1:J1")
Style-HeaderLikeCell $ws.Range("A2:A4")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Lab"
$ws.Range("D2").Value = "Gaussian"
$ws.Range("E2").Value = 0.5075715144200673
$ws.Range("F2").Value = 0.5049922884524501
$ws.Range("G2").Value = 0.5120882634720421
$ws.Range("H2").Value = 0.4975575414068923
$ws.Range("I2").Value = 0.4981614261866363
$ws.Range("J2").Value = 0.5040742067876176

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Lab"
$ws.Range("D3").Value = "Categorical"
$ws.Range("E3").Value = 0.4119384239633809
$ws.Range("F3").Value = 0.4275964632551317
$ws.Range("G3").Value = 0.4409876828343616
$ws.Range("H3").Value = 0.4248359765209714
$ws.Range("I3").Value = 0.4195962607612967
$ws.Range("J3").Value = 0.4249909614670285

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Lab"
$ws.Range("D4").Value = "Bernoulli"
$ws.Range("E4").Value = 0.4119384239633809
$ws.Range("F4").Value = 0.4275964632551317
$ws.Range("G4").Value = 0.4409876828343616
$ws.Range("H4").Value = 0.4248359765209714
$ws.Range("I4").Value = 0.4195962607612967
$ws.Range("J4").Value = 0.4249909614670285

# ---------------------------------------------------------------------------
# Sheet: ClinicalNB3  (identical data to ClinicalNB1)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "ClinicalNB3"

$ws.Range("B1").Value = "Base"
$ws.Range("C1").Value = "P1"
$ws.Range("D1").Value = "P2"
$ws.Range("E1").Value = "P3"
$ws.Range("F1").Value = "P4"
$ws.Range("G1").Value = "P5"
$ws.Range("H1").Value = "Promedio"
$ws.Range("I1").Value = "Accuracy"
Style-HeaderLikeCell $ws.Range("B1:I1")
Style-HeaderLikeCell $ws.Range("A2")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Categorical"
$ws.Range("C2").Value = 0.9216641870017656
$ws.Range("D2").Value = 0.9200612080359398
$ws.Range("E2").Value = 0.9114584800680022
$ws.Range("F2").Value = 0.9086402948228408
$ws.Range("G2").Value = 0.9205627034707913
$ws.Range("H2").Value = 0.916477374679868
$ws.Range("I2").Value = 0.9334

# ---------------------------------------------------------------------------
# Sheet: LabNB1  (Base / P1..P5 / Promedio / Accuracy, 1 data row)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "LabNB1"

$ws.Range("B1").Value = "Base"
$ws.Range("C1").Value = "P1"
$ws.Range("D1").Value = "P2"
$ws.Range("E1").Value = "P3"
$ws.Range("F1").Value = "P4"
$ws.Range("G1").Value = "P5"
$ws.Range("H1").Value = "Promedio"
$ws.Range("I1").Value = "Accuracy"
Style-HeaderLikeCell $ws.Range("B1:I1")
Style-HeaderLikeCell $ws.Range("A2")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Categorical"
$ws.Range("C2").Value = 0.4119384239633809
$ws.Range("D2").Value = 0.4275964632551317
$ws.Range("E2").Value = 0.4409876828343616
$ws.Range("F2").Value = 0.4248359765209714
$ws.Range("G2").Value = 0.4195962607612967
$ws.Range("H2").Value = 0.4249909614670285
$ws.Range("I2").Value = 0.6889999999999999

# ---------------------------------------------------------------------------
# Sheet: LabNB2  (identical data to LabNB1)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "LabNB2"

$ws.Range("B1").Value = "Base"
$ws.Range("C1").Value = "P1"
$ws.Range("D1").Value = "P2"
$ws.Range("E1").Value = "P3"
$ws.Range("F1").Value = "P4"
$ws.Range("G1").Value = "P5"
$ws.Range("H1").Value = "Promedio"
$ws.Range("I1").Value = "Accuracy"
Style-HeaderLikeCell $ws.Range("B1:I1")
Style-HeaderLikeCell $ws.Range("A2")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Categorical"
$ws.Range("C2").Value = 0.4119384239633809
$ws.Range("D2").Value = 0.4275964632551317
$ws.Range("E2").Value = 0.4409876828343616
$ws.Range("F2").Value = 0.4248359765209714
$ws.Range("G2").Value = 0.4195962607612967
$ws.Range("H2").Value = 0.4249909614670285
$ws.Range("I2").Value = 0.6889999999999999

# ---------------------------------------------------------------------------
# Restore original active sheet / selection so the workbook-level view
# state is unaffected by adding the new sheets.
# ---------------------------------------------------------------------------
$origActiveSheet.Activate()
